$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 187
$ws.Range("J2").Value = 198
$ws.Range("L2").Value = 198
$ws.Range("N2").Value = -424
$ws.Range("H51").Value = 2610.875
$ws.Range("I51").Value = 2731.1667
$ws.Range("K51").Value = 2731.1667
$ws.Range("M51").Value = -2247.1667
$ws.Range("H57").Value = 59799.8
$ws.Range("J57").Value = 59499.75
$ws.Range("L57").Value = 178499.25
$ws.Range("N57").Value = -179497.25
$ws.Range("H61").Value = 1115
$ws.Range("I61").Value = 1115
$ws.Range("K61").Value = 3345
$ws.Range("M61").Value = -3173
$ws.Range("H103").Value = 607.5454999999999
$ws.Range("I103").Value = 547.1667
$ws.Range("J103").Value = 680
$ws.Range("K103").Value = 1641.5001
$ws.Range("L103").Value = 2040
$ws.Range("M103").Value = -1055.5001
$ws.Range("N103").Value = -3212
$ws.Range("H112").Value = 1984.091
$ws.Range("J112").Value = 1984.091
$ws.Range("L112").Value = 5952.272999999999
$ws.Range("N112").Value = -8168.272999999999
$ws.Range("H132").Value = 2099.0532
$ws.Range("I132").Value = 1984.8939
$ws.Range("J132").Value = 2936.2222
$ws.Range("K132").Value = 5954.6817
$ws.Range("L132").Value = 8808.6666
$ws.Range("M132").Value = -3424.6817
$ws.Range("N132").Value = -13868.6666
$ws.Range("H133").Value = 49999.742
$ws.Range("J133").Value = 49999.734
$ws.Range("L133").Value = 49999.734
$ws.Range("N133").Value = -60119.734
$ws.Range("H137").Value = 2075.7896
$ws.Range("I137").Value = 2110.3572
$ws.Range("J137").Value = 1979
$ws.Range("K137").Value = 6331.071599999999
$ws.Range("L137").Value = 5937
$ws.Range("M137").Value = -3781.071599999999
$ws.Range("N137").Value = -11037
$ws.Range("H138").Value = 10581.818
$ws.Range("I138").Value = 8799.25
$ws.Range("J138").Value = 10666.702
$ws.Range("K138").Value = 26397.75
$ws.Range("L138").Value = 32000.106
$ws.Range("M138").Value = -21257.75
$ws.Range("N138").Value = -42280.106

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22836.744
$ws.Range("I32").Value = 15988.474
$ws.Range("J32").Value = 28258.291
$ws.Range("K32").Value = 15988.474
$ws.Range("L32").Value = 28258.291
$ws.Range("M32").Value = -15701.474
$ws.Range("N32").Value = -28832.291
$ws.Range("H45").Value = 1471.1177
$ws.Range("J45").Value = 1385.6666
$ws.Range("L45").Value = 1385.6666
$ws.Range("N45").Value = -2139.6666
$ws.Range("H74").Value = 3409.7334
$ws.Range("I74").Value = 2480.4614
$ws.Range("K74").Value = 2480.4614
$ws.Range("M74").Value = -1606.4614
$ws.Range("H77").Value = 3409.7334
$ws.Range("I77").Value = 2480.4614
$ws.Range("K77").Value = 12402.307
$ws.Range("M77").Value = -8034.307000000001
$ws.Range("H132").Value = 4435.078
$ws.Range("I132").Value = 3618.625
$ws.Range("K132").Value = 10855.875
$ws.Range("M132").Value = -8325.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1702800.6
$ws.Range("I86").Value = 2430858.8
$ws.Range("J86").Value = 3998.3333
$ws.Range("K86").Value = 2430858.8
$ws.Range("L86").Value = 3998.3333
$ws.Range("M86").Value = -2429735.8
$ws.Range("N86").Value = -6244.3333
$ws.Range("H89").Value = 1702800.6
$ws.Range("I89").Value = 2430858.8
$ws.Range("J89").Value = 3998.3333
$ws.Range("K89").Value = 12154294
$ws.Range("L89").Value = 19991.6665
$ws.Range("M89").Value = -12148678
$ws.Range("N89").Value = -31223.6665
$ws.Range("H99").Value = 1227.375
$ws.Range("I99").Value = 1188.5
$ws.Range("K99").Value = 1188.5
$ws.Range("M99").Value = 309.5
$ws.Range("H134").Value = 19245.686
$ws.Range("I134").Value = 4559.339
$ws.Range("J134").Value = 98017.91
$ws.Range("K134").Value = 13678.017
$ws.Range("L134").Value = 294053.73
$ws.Range("M134").Value = -11143.017
$ws.Range("N134").Value = -299123.73

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 56355.1
$ws.Range("I31").Value = 4286.467
$ws.Range("J31").Value = 212561
$ws.Range("K31").Value = 4286.467
$ws.Range("L31").Value = 212561
$ws.Range("M31").Value = -3991.467
$ws.Range("N31").Value = -213151
$ws.Range("H34").Value = 56355.1
$ws.Range("I34").Value = 4286.467
$ws.Range("J34").Value = 212561
$ws.Range("K34").Value = 4286.467
$ws.Range("L34").Value = 212561
$ws.Range("M34").Value = -4084.467
$ws.Range("N34").Value = -212965
$ws.Range("H58").Value = 2836.6128
$ws.Range("I58").Value = 2596.4443
$ws.Range("K58").Value = 2596.4443
$ws.Range("M58").Value = -2393.4443
$ws.Range("H99").Value = 5008.5713
$ws.Range("I99").Value = 4247.364
$ws.Range("K99").Value = 4247.364
$ws.Range("M99").Value = -2749.364
$ws.Range("H103").Value = 54999.75
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H126").Value = 5008.5713
$ws.Range("I126").Value = 4247.364
$ws.Range("K126").Value = 12742.092
$ws.Range("M126").Value = -10272.092
$ws.Range("H132").Value = 5113.9
$ws.Range("I132").Value = 3852.2
$ws.Range("J132").Value = 8899
$ws.Range("K132").Value = 11556.6
$ws.Range("L132").Value = 26697
$ws.Range("M132").Value = -9026.599999999999
$ws.Range("N132").Value = -31757
$ws.Range("H134").Value = 387790.3
$ws.Range("I134").Value = 3425.4092
$ws.Range("J134").Value = 2501797.2
$ws.Range("K134").Value = 10276.2276
$ws.Range("L134").Value = 7505391.600000001
$ws.Range("M134").Value = -7741.2276
$ws.Range("N134").Value = -7510461.600000001
$ws.Range("H136").Value = 2836.6128
$ws.Range("I136").Value = 2596.4443
$ws.Range("K136").Value = 7789.3329
$ws.Range("M136").Value = -5239.3329
$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280
$ws.Range("H139").Value = 99402
$ws.Range("J139").Value = 99435
$ws.Range("L139").Value = 99435
$ws.Range("N139").Value = -109715

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 205.42857
$ws.Range("I14").Value = 205.42857
$ws.Range("K14").Value = 616.28571
$ws.Range("M14").Value = -443.28571
$ws.Range("H23").Value = 3283.3333
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15470
$ws.Range("H34").Value = 130072.5
$ws.Range("J34").Value = 207960
$ws.Range("L34").Value = 623880
$ws.Range("N34").Value = -624048
$ws.Range("H133").Value = 33550.83
$ws.Range("I133").Value = 24743.5
$ws.Range("J133").Value = 34960
$ws.Range("K133").Value = 74230.5
$ws.Range("L133").Value = 104880
$ws.Range("M133").Value = -69170.5
$ws.Range("N133").Value = -115000
$ws.Range("H136").Value = 2647.8572
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 71555.64999999999
$ws.Range("I141").Value = 98925.69
$ws.Range("J141").Value = 8995.571
$ws.Range("K141").Value = 296777.07
$ws.Range("L141").Value = 26986.713
$ws.Range("M141").Value = -291597.07
$ws.Range("N141").Value = -37346.713

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 51185.223
$ws.Range("I137").Value = 55349
$ws.Range("J137").Value = 49995.57
$ws.Range("K137").Value = 55349
$ws.Range("L137").Value = 49995.57
$ws.Range("M137").Value = -50249
$ws.Range("N137").Value = -60195.57

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5510.0938
$ws.Range("I132").Value = 4954.3335
$ws.Range("K132").Value = 14863.0005
$ws.Range("M132").Value = -12333.0005
$ws.Range("H136").Value = 11253.885
$ws.Range("I136").Value = 9528.857
$ws.Range("J136").Value = 13266.417
$ws.Range("K136").Value = 28586.571
$ws.Range("L136").Value = 39799.251
$ws.Range("M136").Value = -26036.571
$ws.Range("N136").Value = -44899.251

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 500009250
$ws.Range("J26").Value = 1000000000
$ws.Range("L26").Value = 1000000000
$ws.Range("N26").Value = -1000000586
$ws.Range("H136").Value = 280240.72
$ws.Range("I136").Value = 363677.97
$ws.Range("J136").Value = 134225.5
$ws.Range("K136").Value = 1091033.91
$ws.Range("L136").Value = 402676.5
$ws.Range("M136").Value = -1088483.91
$ws.Range("N136").Value = -407776.5
$ws.Range("H139").Value = 57747.31
$ws.Range("J139").Value = 65357.5
$ws.Range("L139").Value = 65357.5
$ws.Range("N139").Value = -75637.5
